$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2016-11-15 17:24:17"
$wsDeDe.Range("H2").Value = "2016-11-15 17:24:17"
$wsZhCn.Range("H2").Value = "2016-11-15 17:24:02"
$wsZhCn.Range("K2").Value = "2016-11-15 17:24:52"
$wsDeDe.Range("K2").Value = "2016-11-15 17:25:11"
